# Update loading_percent values for the 380 kV case (rows 2-25, columns B-M excluding G)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 17.27400121132155
$ws.Range("C2").Value = 6.210456982042622
$ws.Range("D2").Value = 3.429573696358183
$ws.Range("E2").Value = 10.36743954896778
$ws.Range("F2").Value = 58.69465117415592
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 42.75315948101495
$ws.Range("J2").Value = 10.42460092597672
$ws.Range("K2").Value = 15.86322915611633
$ws.Range("L2").Value = 11.72837058400516
$ws.Range("M2").Value = 17.62854016503052

# Row 3
$ws.Range("B3").Value = 17.2483593757087
$ws.Range("C3").Value = 6.094336099605031
$ws.Range("D3").Value = 3.441385103338726
$ws.Range("E3").Value = 10.38405863763118
$ws.Range("F3").Value = 58.51140703389682
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 42.65336659491704
$ws.Range("J3").Value = 10.43262268719585
$ws.Range("K3").Value = 15.84653044710402
$ws.Range("L3").Value = 11.75000770832539
$ws.Range("M3").Value = 17.65987514556375

# Row 4
$ws.Range("B4").Value = 17.23751755836745
$ws.Range("C4").Value = 6.023962917165669
$ws.Range("D4").Value = 3.449260142005088
$ws.Range("E4").Value = 10.39496118729545
$ws.Range("F4").Value = 58.40547193871022
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 42.59598130303169
$ws.Range("J4").Value = 10.43789362703449
$ws.Range("K4").Value = 15.84037225040131
$ws.Range("L4").Value = 11.76472790485023
$ws.Range("M4").Value = 17.68216696656829

# Row 5
$ws.Range("B5").Value = 17.23433732991727
$ws.Range("C5").Value = 5.995561236524315
$ws.Range("D5").Value = 3.452625987723626
$ws.Range("E5").Value = 10.39958015578754
$ws.Range("F5").Value = 58.36397429323856
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 42.57358143914733
$ws.Range("J5").Value = 10.44012868122948
$ws.Range("K5").Value = 15.83889566593772
$ws.Range("L5").Value = 11.77108772213484
$ws.Range("M5").Value = 17.69201876254461

# Row 6
$ws.Range("B6").Value = 17.23388414955377
$ws.Range("C6").Value = 5.990863088342143
$ws.Range("D6").Value = 3.453194351241379
$ws.Range("E6").Value = 10.40035778165114
$ws.Range("F6").Value = 58.35718510085037
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 42.56992166096047
$ws.Range("J6").Value = 10.44050507737529
$ws.Range("K6").Value = 15.83871294526413
$ws.Range("L6").Value = 11.77216559296541
$ws.Range("M6").Value = 17.69370101781984

# Row 7
$ws.Range("B7").Value = 17.23746965055301
$ws.Range("C7").Value = 6.023578707734196
$ws.Range("D7").Value = 3.449304900350139
$ws.Range("E7").Value = 10.39502276668645
$ws.Range("F7").Value = 58.40490549582012
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 42.5956752127389
$ws.Range("J7").Value = 10.43792341679516
$ws.Range("K7").Value = 15.8403481507212
$ws.Range("L7").Value = 11.76481221250076
$ws.Range("M7").Value = 17.6822967227262

# Row 8
$ws.Range("B8").Value = 17.26414585971584
$ws.Range("C8").Value = 6.170252342280424
$ws.Range("D8").Value = 3.433517069753438
$ws.Range("E8").Value = 10.37302515530144
$ws.Range("F8").Value = 58.63011452670139
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 42.71794860964961
$ws.Range("J8").Value = 10.4272952392535
$ws.Range("K8").Value = 15.85662373624082
$ws.Range("L8").Value = 11.73553348978027
$ws.Range("M8").Value = 17.63871119637766

# Row 9
$ws.Range("B9").Value = 17.35507739569234
$ws.Range("C9").Value = 6.463332278339324
$ws.Range("D9").Value = 3.407494370419141
$ws.Range("E9").Value = 10.33540756791377
$ws.Range("F9").Value = 59.12310140208078
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 42.98822675007965
$ws.Range("J9").Value = 10.40918564381841
$ws.Range("K9").Value = 15.92085323074484
$ws.Range("L9").Value = 11.68948583016178
$ws.Range("M9").Value = 17.57744292963112

# Row 10
$ws.Range("B10").Value = 17.4449801789347
$ws.Range("C10").Value = 6.67949748595061
$ws.Range("D10").Value = 3.391379461507019
$ws.Range("E10").Value = 10.31110540338143
$ws.Range("F10").Value = 59.5154045476521
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 43.2048962935949
$ws.Range("J10").Value = 10.3975331151283
$ws.Range("K10").Value = 15.98745150769691
$ws.Range("L10").Value = 11.66256086518133
$ws.Range("M10").Value = 17.54716116971187

# Row 11
$ws.Range("B11").Value = 17.490774906098
$ws.Range("C11").Value = 6.777515173490927
$ws.Range("D11").Value = 3.384699404859306
$ws.Range("E11").Value = 10.30076777075299
$ws.Range("F11").Value = 59.70012591172854
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 43.30727481102819
$ws.Range("J11").Value = 10.39258816344296
$ws.Range("K11").Value = 16.02188424569288
$ws.Range("L11").Value = 11.65180628615554
$ws.Range("M11").Value = 17.53657695256803

# Row 12
$ws.Range("B12").Value = 17.50880861860662
$ws.Range("C12").Value = 6.814546539023239
$ws.Range("D12").Value = 3.382263314226873
$ws.Range("E12").Value = 10.29695587723512
$ws.Range("F12").Value = 59.77094877830014
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 43.34657993512298
$ws.Range("J12").Value = 10.39076659405812
$ws.Range("K12").Value = 16.03550990821593
$ws.Range("L12").Value = 11.64794813338449
$ws.Range("M12").Value = 17.53302701422131

# Row 13
$ws.Range("B13").Value = 17.50489415074433
$ws.Range("C13").Value = 6.806575608019083
$ws.Range("D13").Value = 3.382783812509389
$ws.Range("E13").Value = 10.29777227410158
$ws.Range("F13").Value = 59.75565743404606
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 43.33809120436423
$ws.Range("J13").Value = 10.39115663748633
$ws.Range("K13").Value = 16.03254942098219
$ws.Range("L13").Value = 11.64876952825333
$ws.Range("M13").Value = 17.53377119819427

# Row 14
$ws.Range("B14").Value = 17.49224474945866
$ws.Range("C14").Value = 6.78056370944168
$ws.Range("D14").Value = 3.384497112692897
$ws.Range("E14").Value = 10.30045210734183
$ws.Range("F14").Value = 59.70593521884729
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 43.31049778346834
$ws.Range("J14").Value = 10.39243728129044
$ws.Range("K14").Value = 16.0229935276169
$ws.Range("L14").Value = 11.65148457974073
$ws.Range("M14").Value = 17.53627572158538

# Row 15
$ws.Range("B15").Value = 17.48458640183351
$ws.Range("C15").Value = 6.764618335021066
$ws.Range("D15").Value = 3.385558733194839
$ws.Range("E15").Value = 10.30210694988452
$ws.Range("F15").Value = 59.6755917685332
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 43.29366553325006
$ws.Range("J15").Value = 10.39322834553202
$ws.Range("K15").Value = 16.01721642165654
$ws.Range("L15").Value = 11.65317553198787
$ws.Range("M15").Value = 17.53786944319823

# Row 16
$ws.Range("B16").Value = 17.44208493367594
$ws.Range("C16").Value = 6.673081939490817
$ws.Range("D16").Value = 3.391829106863296
$ws.Range("E16").Value = 10.31179539324971
$ws.Range("F16").Value = 59.50345624164369
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 43.19828136973982
$ws.Range("J16").Value = 10.39786342326212
$ws.Range("K16").Value = 15.9852838134668
$ws.Range("L16").Value = 11.66329372403579
$ws.Range("M16").Value = 17.54791702097949

# Row 17
$ws.Range("B17").Value = 17.41725782195046
$ws.Range("C17").Value = 6.616816928572092
$ws.Range("D17").Value = 3.3958423894432
$ws.Range("E17").Value = 10.31792240659411
$ws.Range("F17").Value = 59.39944111234872
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 43.14073497908582
$ws.Range("J17").Value = 10.40079789517534
$ws.Range("K17").Value = 15.96674819989636
$ws.Range("L17").Value = 11.66988319060018
$ws.Range("M17").Value = 17.55489767630952

# Row 18
$ws.Range("B18").Value = 17.4034396880505
$ws.Range("C18").Value = 6.584427242369308
$ws.Range("D18").Value = 3.398211974221178
$ws.Range("E18").Value = 10.3215140672376
$ws.Range("F18").Value = 59.34020536981552
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 43.10799560615986
$ws.Range("J18").Value = 10.40251923081239
$ws.Range("K18").Value = 15.95647692356764
$ws.Range("L18").Value = 11.67381390015901
$ws.Range("M18").Value = 17.55921318333104

# Row 19
$ws.Range("B19").Value = 17.39884076222009
$ws.Range("C19").Value = 6.573457177824367
$ws.Range("D19").Value = 3.399024796810997
$ws.Range("E19").Value = 10.32274175951849
$ws.Range("F19").Value = 59.32025148493758
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 43.09697276356895
$ws.Range("J19").Value = 10.40310780590636
$ws.Range("K19").Value = 15.95306644975152
$ws.Range("L19").Value = 11.67516893553048
$ws.Range("M19").Value = 17.56072595927922

# Row 20
$ws.Range("B20").Value = 17.41985300312341
$ws.Range("C20").Value = 6.622809582802653
$ws.Range("D20").Value = 3.395408829928685
$ws.Range("E20").Value = 10.31726318607172
$ws.Range("F20").Value = 59.41045272742952
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 43.14682373796516
$ws.Range("J20").Value = 10.40048204939122
$ws.Range("K20").Value = 15.96868104894801
$ws.Range("L20").Value = 11.66916717971479
$ws.Range("M20").Value = 17.55412348678646

# Row 21
$ws.Range("B21").Value = 17.49594149893275
$ws.Range("C21").Value = 6.788206674344598
$ws.Range("D21").Value = 3.38399133783211
$ws.Range("E21").Value = 10.29966219040261
$ws.Range("F21").Value = 59.72051637314711
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 43.31858817496087
$ws.Range("J21").Value = 10.39205974320742
$ws.Range("K21").Value = 16.02578447020932
$ws.Range("L21").Value = 11.65068128874757
$ws.Range("M21").Value = 17.53552765740852

# Row 22
$ws.Range("B22").Value = 17.54969834477552
$ws.Range("C22").Value = 6.895785833737525
$ws.Range("D22").Value = 3.377074314416327
$ws.Range("E22").Value = 10.28875758978376
$ws.Range("F22").Value = 59.92823520880135
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 43.43396791521268
$ws.Range("J22").Value = 10.38685232054133
$ws.Range("K22").Value = 16.06652045058882
$ws.Range("L22").Value = 11.63984901017574
$ws.Range("M22").Value = 17.52604374948139

# Row 23
$ws.Range("B23").Value = 17.52064292510169
$ws.Range("C23").Value = 6.838428937066071
$ws.Range("D23").Value = 3.380716218158277
$ws.Range("E23").Value = 10.29452294610883
$ws.Range("F23").Value = 59.81691682760221
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 43.37210603590744
$ws.Range("J23").Value = 10.38960450421795
$ws.Range("K23").Value = 16.04446925705941
$ws.Range("L23").Value = 11.64551623279456
$ws.Range("M23").Value = 17.53086152720588

# Row 24
$ws.Range("B24").Value = 17.41867830255293
$ws.Range("C24").Value = 6.620100432135024
$ws.Range("D24").Value = 3.395604648155012
$ws.Range("E24").Value = 10.31756100422241
$ws.Range("F24").Value = 59.40547261734839
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 43.14406993496139
$ws.Range("J24").Value = 10.40062473651132
$ws.Range("K24").Value = 15.96780600771604
$ws.Range("L24").Value = 11.6694904448303
$ws.Range("M24").Value = 17.55447255634995

# Row 25
$ws.Range("B25").Value = 17.32638519142507
$ws.Range("C25").Value = 6.383723089846399
$ws.Range("D25").Value = 3.414006243167162
$ws.Range("E25").Value = 10.34499626285973
$ws.Range("F25").Value = 58.98435577580348
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 42.91189267283754
$ws.Range("J25").Value = 10.41379361887622
$ws.Range("K25").Value = 15.90004473747404
$ws.Range("L25").Value = 11.70072834054314
$ws.Range("M25").Value = 17.591428387428
